# Rename orgz module to party and use uuid for FK rels
#
# 1. Rename the "organisations" sheet to "organizations" (American spelling).
# 2. On that sheet, the "identifier" column header becomes "code" and the
#    "fncode" column header becomes "subtype" (the underlying data values in
#    those two columns are unchanged - only the header labels change).
# 3. Update the sheet's current selection / view position.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("organisations")
$ws.Name = "organizations"

# Header row (row 1): column B "identifier" -> "code", column C "fncode" -> "subtype"
$ws.Cells.Item(1, 2).Value = "code"
$ws.Cells.Item(1, 3).Value = "subtype"

# Move the active selection/view to reflect the edited sheet state
$ws.Activate()
$ws.Range("T3").Select()
